$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeout) values for rows 2-36, replacing previous Strike# values.
$newK = @{
    2  = 7
    3  = 5
    4  = 6
    5  = 5
    6  = 2
    7  = 7
    8  = 6
    9  = 5
    10 = 2
    11 = 4
    12 = 6
    13 = 1
    14 = 14
    15 = 6
    16 = 3
    17 = 5
    18 = 5
    19 = 5
    20 = 6
    21 = 5
    22 = 5
    23 = 6
    24 = 5
    25 = 8
    26 = 2
    27 = 11
    28 = 5
    29 = 7
    30 = 6
    31 = 9
    32 = 2
    33 = 6
    34 = 4
    35 = 3
    36 = 4
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
